# BioFoundries.xlsx - "Add files via upload" - append two new BioFoundry
# entries (VTT / Finland and SynBio Foundry SJTU / Shanghai) to Sheet1.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 30: VTT Technical Research Centre of Finland Ltd -----------------
$ws.Range("A30").Value = "VTT Technical Research Centre of Finland Ltd"
$ws.Range("B30").Value = "Espoo"
$ws.Range("C30").Value = "Finland"
$ws.Range("D30").Value = "FIN"
$ws.Range("E30").Value = 23
$ws.Range("F30").Value = "FIN"
$ws.Range("G30").Value = 60.1
$ws.Range("H30").Value = 24.4

# --- Row 31: SynBio Foundry SJTU, Shanghai Jiao Tong University -----------
$ws.Range("A31").Value = "SynBio Foundry SJTU, Shanghai Jiao Tong University"
$ws.Range("B31").Value = "Shanghai"
$ws.Range("C31").Value = "People's Republic of China"
$ws.Range("D31").Value = "SJT"
$ws.Range("E31").Value = 24
$ws.Range("F31").Value = "CHN"
$ws.Range("G31").Value = 31.1
$ws.Range("H31").Value = 121.3

# The name cell for the new Shanghai entry was (re)typed with a slightly
# larger font than the sheet default, which also grows that row's height.
$ws.Range("A31").Font.Size = 12
$ws.Range("A31:H31").RowHeight = 15.5

# Leave the selection where the author ended up after entering the data.
$ws.Range("F32").Select()
